$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-12 Sunday", "2024-05-13 Monday"),
    @("575÷5=115, 0", "594÷7=84, 6"),
    @("442÷3=147, 1", "750÷2=375, 0"),
    @("869÷6=144, 5", "799÷6=133, 1"),
    @("491÷3=163, 2", "554÷4=138, 2"),
    @("908÷2=454, 0", "725÷7=103, 4"),
    @("679÷4=169, 3", "988÷7=141, 1"),
    @("148÷9=16, 4", "234÷2=117, 0"),
    @("180÷8=22, 4", "440÷3=146, 2"),
    @("627÷7=89, 4", "823÷6=137, 1"),
    @("227÷4=56, 3", "728÷4=182, 0"),
    @("458÷2=229, 0", "571÷2=285, 1"),
    @("894÷9=99, 3", "420÷4=105, 0"),
    @("915÷5=183, 0", "746÷7=106, 4"),
    @("887÷9=98, 5", "671÷9=74, 5"),
    @("675÷5=135, 0", "720÷9=80, 0"),
    @("782÷3=260, 2", "788÷5=157, 3"),
    @("270÷5=54, 0", "867÷9=96, 3"),
    @("759÷4=189, 3", "963÷3=321, 0"),
    @("742÷5=148, 2", "753÷8=94, 1"),
    @("714÷8=89, 2", "793÷9=88, 1"),
    @("373÷2=186, 1", "848÷8=106, 0"),
    @("885÷3=295, 0", "542÷5=108, 2"),
    @("124÷6=20, 4", "308÷9=34, 2"),
    @("605÷3=201, 2", "668÷5=133, 3"),
    @("168÷7=24, 0", "175÷6=29, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
